# keyWords.xlsx — remove the "output" entry from the keyword list (column B)
# and shift the remaining entries up by one row, the way Excel behaves when
# a single cell is deleted with "Shift cells up" inside a column.
#
# Before (col B, rows 3-9): output, control, handle, dynamic, update, do, controller
# After  (col B, rows 3-8):          control, handle, dynamic, update, do, controller

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift B4:B9 into B3:B8 (column B only — the other columns are untouched).
for ($r = 3; $r -le 8; $r++) {
    $ws.Range("B$r").Value = $ws.Range("B" + ($r + 1)).Value()
}

# The last row's old value has now moved up, so clear the vacated cell.
$ws.Range("B9").ClearContents()

# The saved file shows the selection sitting on B3.
[void]$ws.Range("B3").Select()

# Match the recorded window geometry at save time.
$win = $wb.Windows.Item(1)
$win.Width = 21600
$win.Height = 10695
$win.Left = 29940
$win.Top = 3600
